$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual student attendance statuses
$ws.Range("D2").Value = "Hadir"
$ws.Range("D18").Value = "Sakit"
$ws.Range("D24").Value = "Izin"
$ws.Range("D25").Value = "Izin"
$ws.Range("D34").Value = "Alpha"

# Update the summary counts
$ws.Range("A37").Value = "Hadir: 29"
$ws.Range("A38").Value = "Izin: 2"
$ws.Range("A39").Value = "Sakit: 1"
